$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LQFP32")

# Fill in the previously-empty "uC IO function" (D) column for the
# GPIO-based pin rows (EN pins -> generic GPIO use).
$ws.Range("D12").Value = "GPIO"
$ws.Range("D13").Value = "GPIO"
$ws.Range("D14").Value = "GPIO"

# Remove the four "high side switch enable" (AUX) rows - rows 15-18 -
# shifting the EXTI0/1/2 rows up in their place.
$ws.Range("A15:I18").EntireRow.Delete() | Out-Null

# Now rows 19,20,21 have become rows 15,16,17 - fill in their
# previously-empty "uC IO function" (D) column.
$ws.Range("D15").Value = "EXTI0"
$ws.Range("D16").Value = "EXTI1"
$ws.Range("D17").Value = "EXTI2"

# Keep the autofilter range in sync with the new (smaller) table extent
# (matches the saved filter range, which kept its 5-row slack below the data).
$ws.AutoFilterMode = $false
$ws.Range("A1:I22").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name tracks the filter range too.
$wb.Names.Item("LQFP32!_FilterDatabase").RefersTo = "=LQFP32!`$A`$1:`$I`$22"

# Move the selection/active cell (matches the saved UI state in the diff).
$ws.Range("D14").Select() | Out-Null
